$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 3.3
$ws.Range("I2").Value = 2.4
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 1.95
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("AA2").Value = 2
$ws.Range("AB2").Value = 1.75
$ws.Range("AI2").Value = 7
$ws.Range("AK2").Value = 17
$ws.Range("AN2").Value = 6.5
$ws.Range("AO2").Value = 10

# Row 4 updates
$ws.Range("I4").Value = 4.1
$ws.Range("N4").Value = 8.5
$ws.Range("AB4").Value = 1.73
$ws.Range("AF4").Value = 15

# Row 5 updates
$ws.Range("G5").Value = 1.57
$ws.Range("I5").Value = 5.5
$ws.Range("AB5").Value = 1.67
$ws.Range("AJ5").Value = 7
$ws.Range("AM5").Value = 1250
